$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 1: @@ -7298,25 +7298,25 @@
$ws.Range("H132").Value = 2720709
$ws.Range("I132").Value = 2910264.8
$ws.Range("J132").Value = 3743.3333
$ws.Range("K132").Value = 8730794.399999999
$ws.Range("L132").Value = 11229.9999
$ws.Range("M132").Value = -8728264.399999999
$ws.Range("N132").Value = -16289.9999

$ws = $wb.Worksheets.Item("ARM")
# Hunk 2: @@ -7897,22 +7897,22 @@
$ws.Range("H2").Value = 1051.375
$ws.Range("I2").Value = 1051.375
$ws.Range("K2").Value = 1051.375
$ws.Range("M2").Value = -938.375

# Hunk 3: @@ -10037,22 +10037,22 @@
$ws.Range("H45").Value = 1820
$ws.Range("I45").Value = 1820
$ws.Range("K45").Value = 1820
$ws.Range("M45").Value = -1443

# Hunk 4: @@ -10934,19 +10934,25 @@
$ws.Range("H63").Value = 2333
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 999
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 999
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -2371

# Hunk 5: @@ -11078,19 +11084,25 @@
$ws.Range("H66").Value = 2333
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 999
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 4995
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -11859

# Hunk 6: @@ -12147,25 +12159,25 @@
$ws.Range("H88").Value = 2683.8572
$ws.Range("I88").Value = 1493
$ws.Range("J88").Value = 3160.2
$ws.Range("K88").Value = 1493
$ws.Range("L88").Value = 3160.2
$ws.Range("M88").Value = -1087
$ws.Range("N88").Value = -3972.2

# Hunk 7: @@ -12297,25 +12309,25 @@
$ws.Range("H91").Value = 2683.8572
$ws.Range("I91").Value = 1493
$ws.Range("J91").Value = 3160.2
$ws.Range("K91").Value = 1493
$ws.Range("L91").Value = 3160.2
$ws.Range("M91").Value = -89
$ws.Range("N91").Value = -5968.2

# Hunk 8: @@ -12839,25 +12851,25 @@
$ws.Range("H102").Value = 3025.32
$ws.Range("I102").Value = 2723.3157
$ws.Range("J102").Value = 3981.6667
$ws.Range("K102").Value = 2723.3157
$ws.Range("L102").Value = 3981.6667
$ws.Range("M102").Value = -1101.3157
$ws.Range("N102").Value = -7225.6667

# Hunk 9: @@ -13513,22 +13525,22 @@
$ws.Range("H116").Value = 1051.375
$ws.Range("I116").Value = 1051.375
$ws.Range("K116").Value = 1051.375
$ws.Range("M116").Value = 1242.625

$ws = $wb.Worksheets.Item("BSM")
# Hunk 10: @@ -14918,22 +14930,22 @@
$ws.Range("H3").Value = 1051.375
$ws.Range("I3").Value = 1051.375
$ws.Range("K3").Value = 1051.375
$ws.Range("M3").Value = -937.375

# Hunk 11: @@ -18771,22 +18783,22 @@
$ws.Range("H82").Value = 10047.167
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2617

# Hunk 12: @@ -18924,22 +18936,22 @@
$ws.Range("H85").Value = 10047.167
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1674

# Hunk 13: @@ -19619,25 +19631,25 @@
$ws.Range("H99").Value = 1585.5625
$ws.Range("I99").Value = 1550.75
$ws.Range("J99").Value = 1690
$ws.Range("K99").Value = 1550.75
$ws.Range("L99").Value = 1690
$ws.Range("M99").Value = -52.75
$ws.Range("N99").Value = -4686

# Hunk 14: @@ -19910,25 +19922,25 @@
$ws.Range("H105").Value = 5513.3335
$ws.Range("I105").Value = 7145
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 7145
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = -5398
$ws.Range("N105").Value = -5744

$ws = $wb.Worksheets.Item("CRP")
# Hunk 15: @@ -28342,25 +28354,25 @@
$ws.Range("H134").Value = 1379.2693
$ws.Range("I134").Value = 1332.6316
$ws.Range("J134").Value = 1505.8572
$ws.Range("K134").Value = 3997.8948
$ws.Range("L134").Value = 4517.571599999999
$ws.Range("M134").Value = -1462.8948
$ws.Range("N134").Value = -9587.571599999999

$ws = $wb.Worksheets.Item("CUL")
# Hunk 16: @@ -34456,25 +34468,25 @@
$ws.Range("H113").Value = 1165.3478
$ws.Range("I113").Value = 1823.1818
$ws.Range("J113").Value = 562.3333
$ws.Range("K113").Value = 5469.5454
$ws.Range("L113").Value = 1686.9999
$ws.Range("M113").Value = -3299.5454
$ws.Range("N113").Value = -6026.9999

# Hunk 17: @@ -35432,25 +35444,25 @@
$ws.Range("H132").Value = 1529
$ws.Range("I132").Value = 1384
$ws.Range("J132").Value = 1587
$ws.Range("K132").Value = 12456
$ws.Range("L132").Value = 14283
$ws.Range("M132").Value = -9926
$ws.Range("N132").Value = -19343

# Hunk 18: @@ -35692,25 +35704,25 @@
$ws.Range("H137").Value = 2887.3157
$ws.Range("I137").Value = 1188
$ws.Range("J137").Value = 5800.4287
$ws.Range("K137").Value = 3564
$ws.Range("L137").Value = 17401.2861
$ws.Range("M137").Value = 1536
$ws.Range("N137").Value = -27601.2861

$ws = $wb.Worksheets.Item("GSM")
# Hunk 19: @@ -37151,22 +37163,22 @@
$ws.Range("H24").Value = 33756.125
$ws.Range("J24").Value = 33756.125
$ws.Range("L24").Value = 33756.125
$ws.Range("N24").Value = -34102.125

# Hunk 20: @@ -39390,22 +39402,22 @@
$ws.Range("H70").Value = 5157
$ws.Range("I70").Value = 4945.684
$ws.Range("K70").Value = 4945.684
$ws.Range("M70").Value = -4675.684

# Hunk 21: @@ -39537,22 +39549,22 @@
$ws.Range("H73").Value = 5157
$ws.Range("I73").Value = 4945.684
$ws.Range("K73").Value = 4945.684
$ws.Range("M73").Value = -4009.684

# Hunk 22: @@ -39871,25 +39883,19 @@
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

# Hunk 23: @@ -40021,25 +40027,19 @@
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Hunk 24: @@ -46935,25 +46935,25 @@
$ws.Range("H82").Value = 1880.7715
$ws.Range("I82").Value = 825.25
$ws.Range("J82").Value = 2193.5186
$ws.Range("K82").Value = 825.25
$ws.Range("L82").Value = 2193.5186
$ws.Range("M82").Value = -464.25
$ws.Range("N82").Value = -2915.5186

# Hunk 25: @@ -47079,25 +47079,25 @@
$ws.Range("H85").Value = 1880.7715
$ws.Range("I85").Value = 825.25
$ws.Range("J85").Value = 2193.5186
$ws.Range("K85").Value = 825.25
$ws.Range("L85").Value = 2193.5186
$ws.Range("M85").Value = 422.75
$ws.Range("N85").Value = -4689.518599999999

$ws = $wb.Worksheets.Item("WVR")
# Hunk 26: @@ -55892,22 +55892,22 @@
$ws.Range("H123").Value = 26000
$ws.Range("J123").Value = 26000
$ws.Range("L123").Value = 26000
$ws.Range("N123").Value = -35800

# Hunk 27: @@ -56336,25 +56336,25 @@
$ws.Range("H132").Value = 7713
$ws.Range("I132").Value = 8750.134
$ws.Range("J132").Value = 4601.6
$ws.Range("K132").Value = 26250.402
$ws.Range("L132").Value = 13804.8
$ws.Range("M132").Value = -23720.402
$ws.Range("N132").Value = -18864.8

# Hunk 28: @@ -56538,25 +56538,25 @@
$ws.Range("H136").Value = 722.1667
$ws.Range("I136").Value = 648.05884
$ws.Range("J136").Value = 902.1429000000001
$ws.Range("K136").Value = 1944.17652
$ws.Range("L136").Value = 2706.4287
$ws.Range("M136").Value = 605.82348
$ws.Range("N136").Value = -7806.4287
